# Add cantrals by cantons
# Reshape the header from a two-row, merged-looking header into a single
# header row (idx, idx2, Name, Date Start, Date End, (m3/s), (MW1), (MW2),
# (GWh) Winter, (GWh) Summer, (GWh) Year) immediately above the data, and
# drop the trailing blank/old header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old sheet had two header rows (row 1 with scattered "mation"/"pompes)"/
# "Hiver"/"Eté"/"Année" fragments, row 2 with unit labels) followed by 13
# rows of data (rows 3-15). Delete the first (top) header row entirely -
# this shifts everything up by one, turning the old unit-label row into the
# new row 1, and the 13 data rows into rows 2-14.
$ws.Rows.Item(1).Delete()

# Rewrite row 1 as the single consolidated header row.
$ws.Cells.Item(1, 1).Value = "idx"
$ws.Cells.Item(1, 2).Value = "idx2"
$ws.Cells.Item(1, 3).Value = "Name"
$ws.Cells.Item(1, 4).Value = "Date Start"
$ws.Cells.Item(1, 5).Value = "Date End"
$ws.Cells.Item(1, 6).Value = "(m3/s)"
$ws.Cells.Item(1, 7).Value = "(MW1)"
$ws.Cells.Item(1, 8).Value = "(MW2)"
$ws.Cells.Item(1, 9).Value = "(GWh) Winter"
$ws.Cells.Item(1, 10).Value = "(GWh) Summer"
$ws.Cells.Item(1, 11).Value = "(GWh) Year"

# Columns A-E of the header use the plain default font/format.
# Columns F-K keep the small (9pt Arial) header font used elsewhere in the
# sheet, with the default "General" number format.
$ws.Range("F1:K1").Font.Name = "Arial"
$ws.Range("F1:K1").Font.Size = 9

$ws.Range("A2:K2").Select()
